# Impulse rev C BOM update: remove the obsolete "0k" / R17 resistor line
# (row 25). Excel will shift the rows below it (U1..X2 parts) up by one,
# renumbering rows 26-31 to 25-30 and compacting the shared strings table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Delete()

# Match the author's final cursor position recorded in the saved file.
$ws.Range("E16").Select()
